$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# Row 12: add an end time (C12, 17:30 as an Excel time fraction), which updates
# the worked-time formula in D12, set the row height to fit the new task
# description, and fill in the task text (F12).
$ws.Range("C12").Value = 0.72916666666666663
$ws.Rows.Item(12).RowHeight = 60
$ws.Range("F12").Value = "30 min ohjauskokous, 1h pöytäkirjan teko, 2h Projektin tiedostojen siirto OneDrive:stä GitHub:iin sekä OneNote projekti taulukkojen muokkaaminen md-tiedostomuotoon GitHub:ssa. 5 min WordPress blogikirjoitusta."

# Update the view to reflect where the author was working when saving.
$ws.Activate()
$ws.Range("F13").Select()
